# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the leve-profit sheets to match the refreshed market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 8568
$ws.Range("I13").Value = 7000
$ws.Range("J13").Value = 8960
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 8960
$ws.Range("M13").Value = -6831
$ws.Range("N13").Value = -9298


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1683.8928
$ws.Range("I45").Value = 1774.3846
$ws.Range("K45").Value = 1774.3846
$ws.Range("M45").Value = -1397.3846

$ws.Range("H46").Value = 12129.125
$ws.Range("I46").Value = 4329.3335
$ws.Range("J46").Value = 16809
$ws.Range("K46").Value = 4329.3335
$ws.Range("L46").Value = 16809
$ws.Range("M46").Value = -4010.3335
$ws.Range("N46").Value = -17447

$ws.Range("H61").Value = 5865.0938
$ws.Range("I61").Value = 3866.5908
$ws.Range("J61").Value = 10261.8
$ws.Range("K61").Value = 3866.5908
$ws.Range("L61").Value = 10261.8
$ws.Range("M61").Value = -3654.5908
$ws.Range("N61").Value = -10685.8

$ws.Range("H74").Value = 2842.422
$ws.Range("I74").Value = 2129.25
$ws.Range("J74").Value = 4017.0588
$ws.Range("K74").Value = 2129.25
$ws.Range("L74").Value = 4017.0588
$ws.Range("M74").Value = -1255.25
$ws.Range("N74").Value = -5765.0588

$ws.Range("H77").Value = 2842.422
$ws.Range("I77").Value = 2129.25
$ws.Range("J77").Value = 4017.0588
$ws.Range("K77").Value = 10646.25
$ws.Range("L77").Value = 20085.294
$ws.Range("M77").Value = -6278.25
$ws.Range("N77").Value = -28821.294

$ws.Range("H110").Value = 1632.0834
$ws.Range("I110").Value = 1430.303
$ws.Range("J110").Value = 3851.6667
$ws.Range("K110").Value = 1430.303
$ws.Range("L110").Value = 3851.6667
$ws.Range("M110").Value = 614.6969999999999
$ws.Range("N110").Value = -7941.6667

$ws.Range("H122").Value = 2140.5217
$ws.Range("I122").Value = 1719.1
$ws.Range("K122").Value = 5157.299999999999
$ws.Range("M122").Value = -2707.299999999999

$ws.Range("H132").Value = 3581.7778
$ws.Range("I132").Value = 3598.3428
$ws.Range("J132").Value = 3002
$ws.Range("K132").Value = 10795.0284
$ws.Range("L132").Value = 9006
$ws.Range("M132").Value = -8265.028399999999
$ws.Range("N132").Value = -14066

$ws.Range("H136").Value = 5865.0938
$ws.Range("I136").Value = 3866.5908
$ws.Range("J136").Value = 10261.8
$ws.Range("K136").Value = 11599.7724
$ws.Range("L136").Value = 30785.4
$ws.Range("M136").Value = -9049.7724
$ws.Range("N136").Value = -35885.39999999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 36072.75
$ws.Range("J44").Value = 36072.75
$ws.Range("L44").Value = 36072.75
$ws.Range("N44").Value = -37066.75

$ws.Range("H134").Value = 8282.217000000001
$ws.Range("I134").Value = 4520.82
$ws.Range("J134").Value = 27089.2
$ws.Range("K134").Value = 13562.46
$ws.Range("L134").Value = 81267.60000000001
$ws.Range("M134").Value = -11027.46
$ws.Range("N134").Value = -86337.60000000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2942.0222
$ws.Range("I31").Value = 1284.2222
$ws.Range("J31").Value = 3356.4722
$ws.Range("K31").Value = 1284.2222
$ws.Range("L31").Value = 3356.4722
$ws.Range("M31").Value = -989.2221999999999
$ws.Range("N31").Value = -3946.4722

$ws.Range("H34").Value = 2942.0222
$ws.Range("I34").Value = 1284.2222
$ws.Range("J34").Value = 3356.4722
$ws.Range("K34").Value = 1284.2222
$ws.Range("L34").Value = 3356.4722
$ws.Range("M34").Value = -1082.2222
$ws.Range("N34").Value = -3760.4722

$ws.Range("H107").Value = 1498.2727
$ws.Range("I107").Value = 992.5
$ws.Range("K107").Value = 992.5
$ws.Range("M107").Value = 927.5

$ws.Range("H122").Value = 2452.0715
$ws.Range("I122").Value = 2071.1428
$ws.Range("K122").Value = 6213.428400000001
$ws.Range("M122").Value = -3763.428400000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1113.5
$ws.Range("I5").Value = 788.5806
$ws.Range("J5").Value = 2232.6667
$ws.Range("K5").Value = 2365.7418
$ws.Range("L5").Value = 6698.000100000001
$ws.Range("M5").Value = -2253.7418
$ws.Range("N5").Value = -6922.000100000001

$ws.Range("H12").Value = 983.63635
$ws.Range("J12").Value = 1195.5555
$ws.Range("L12").Value = 3586.6665
$ws.Range("N12").Value = -3932.6665

$ws.Range("H69").Value = 7275.4443
$ws.Range("J69").Value = 7559.875
$ws.Range("L69").Value = 22679.625
$ws.Range("N69").Value = -24301.625

$ws.Range("H72").Value = 7275.4443
$ws.Range("J72").Value = 7559.875
$ws.Range("L72").Value = 68038.875
$ws.Range("N72").Value = -76150.875

$ws.Range("H122").Value = 10001376
$ws.Range("J122").Value = 14286826
$ws.Range("L122").Value = 128581434
$ws.Range("N122").Value = -128586334

$ws.Range("H126").Value = 5999.5
$ws.Range("I126").Value = 5999.5
$ws.Range("K126").Value = 17998.5
$ws.Range("M126").Value = -13058.5

$ws.Range("H135").Value = 1113.5
$ws.Range("I135").Value = 788.5806
$ws.Range("J135").Value = 2232.6667
$ws.Range("K135").Value = 7097.2254
$ws.Range("L135").Value = 20094.0003
$ws.Range("M135").Value = -4562.2254
$ws.Range("N135").Value = -25164.0003


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 35499.1
$ws.Range("J57").Value = 47147.285
$ws.Range("L57").Value = 47147.285
$ws.Range("N57").Value = -48787.285

$ws.Range("H80").Value = 7635.391
$ws.Range("I80").Value = 6513.75
$ws.Range("J80").Value = 10199.143
$ws.Range("K80").Value = 6513.75
$ws.Range("L80").Value = 10199.143
$ws.Range("M80").Value = -5515.75
$ws.Range("N80").Value = -12195.143

$ws.Range("H83").Value = 7635.391
$ws.Range("I83").Value = 6513.75
$ws.Range("J83").Value = 10199.143
$ws.Range("K83").Value = 32568.75
$ws.Range("L83").Value = 50995.715
$ws.Range("M83").Value = -27576.75
$ws.Range("N83").Value = -60979.715

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 183680.55
$ws.Range("I113").Value = 1998.625
$ws.Range("J113").Value = 668165.7
$ws.Range("K113").Value = 1998.625
$ws.Range("L113").Value = 668165.7
$ws.Range("M113").Value = 171.375
$ws.Range("N113").Value = -672505.7

$ws.Range("H132").Value = 10948.25
$ws.Range("I132").Value = 13440.579
$ws.Range("K132").Value = 40321.737
$ws.Range("M132").Value = -37791.737


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1453.9615
$ws.Range("J46").Value = 1592.7
$ws.Range("L46").Value = 1592.7
$ws.Range("N46").Value = -1968.7

$ws.Range("H61").Value = 1092
$ws.Range("I61").Value = 1092
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1092
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -890
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 2716.653
$ws.Range("I68").Value = 2427.3171
$ws.Range("J68").Value = 4199.5
$ws.Range("K68").Value = 2427.3171
$ws.Range("L68").Value = 4199.5
$ws.Range("M68").Value = -1678.3171
$ws.Range("N68").Value = -5697.5

$ws.Range("H71").Value = 2716.653
$ws.Range("I71").Value = 2427.3171
$ws.Range("J71").Value = 4199.5
$ws.Range("K71").Value = 12136.5855
$ws.Range("L71").Value = 20997.5
$ws.Range("M71").Value = -8392.585500000001
$ws.Range("N71").Value = -28485.5

$ws.Range("H82").Value = 2047.238
$ws.Range("J82").Value = 2648.4443
$ws.Range("L82").Value = 2648.4443
$ws.Range("N82").Value = -3370.4443

$ws.Range("H85").Value = 2047.238
$ws.Range("J85").Value = 2648.4443
$ws.Range("L85").Value = 2648.4443
$ws.Range("N85").Value = -5144.4443

$ws.Range("H113").Value = 1092
$ws.Range("I113").Value = 1092
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1092
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1078
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 3537.9768
$ws.Range("I132").Value = 3133.2258
$ws.Range("K132").Value = 9399.6774
$ws.Range("M132").Value = -6869.6774

$ws.Range("H136").Value = 2357.4426
$ws.Range("I136").Value = 2059
$ws.Range("K136").Value = 6177
$ws.Range("M136").Value = -3627


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 3599
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H62").Value = 4209.1904
$ws.Range("I62").Value = 3745.3845
$ws.Range("K62").Value = 3745.3845
$ws.Range("M62").Value = -3121.3845

$ws.Range("H65").Value = 4209.1904
$ws.Range("I65").Value = 3745.3845
$ws.Range("K65").Value = 18726.9225
$ws.Range("M65").Value = -15606.9225

$ws.Range("H107").Value = 1168.4565
$ws.Range("I107").Value = 1015
$ws.Range("J107").Value = 1407.1666
$ws.Range("K107").Value = 3045
$ws.Range("L107").Value = 4221.4998
$ws.Range("M107").Value = -1125
$ws.Range("N107").Value = -8061.4998

$ws.Range("H113").Value = 581.1875
$ws.Range("I113").Value = 432.23077
$ws.Range("K113").Value = 1296.69231
$ws.Range("M113").Value = 873.3076900000001

$ws.Range("H136").Value = 751.13336
$ws.Range("I136").Value = 756.6279
$ws.Range("K136").Value = 2269.8837
$ws.Range("M136").Value = 280.1163000000001

